# Appends the 10 newly-played NBA games (row 718-727) to Sheet1,
# mirroring the football/basketball results table layout already
# used by the sheet: Away team, Away Pts, Home team, Home Pts,
# Overtime, Attend., Arena, Win, Loss.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(718,1).Value = "Los Angeles Clippers"
$ws.Cells.Item(718,2).Value = 136
$ws.Cells.Item(718,2).NumberFormat = "#,##0"
$ws.Cells.Item(718,3).Value = "Detroit Pistons"
$ws.Cells.Item(718,4).Value = 125
$ws.Cells.Item(718,4).NumberFormat = "#,##0"
$ws.Cells.Item(718,5).Value = "No"
$ws.Cells.Item(718,6).Value = 17832
$ws.Cells.Item(718,7).Value = "Little Caesars Arena"
$ws.Cells.Item(718,8).Value = "Los Angeles Clippers"
$ws.Cells.Item(718,9).Value = "Detroit Pistons"

$ws.Cells.Item(719,1).Value = "Miami Heat"
$ws.Cells.Item(719,2).Value = 110
$ws.Cells.Item(719,2).NumberFormat = "#,##0"
$ws.Cells.Item(719,3).Value = "Washington Wizards"
$ws.Cells.Item(719,4).Value = 102
$ws.Cells.Item(719,4).NumberFormat = "#,##0"
$ws.Cells.Item(719,5).Value = "No"
$ws.Cells.Item(719,6).Value = 17832
$ws.Cells.Item(719,7).Value = "Capital One Arena"
$ws.Cells.Item(719,8).Value = "Miami Heat"
$ws.Cells.Item(719,9).Value = "Washington Wizards"

$ws.Cells.Item(720,1).Value = "Phoenix Suns"
$ws.Cells.Item(720,2).Value = 120
$ws.Cells.Item(720,2).NumberFormat = "#,##0"
$ws.Cells.Item(720,3).Value = "Atlanta Hawks"
$ws.Cells.Item(720,4).Value = 129
$ws.Cells.Item(720,4).NumberFormat = "#,##0"
$ws.Cells.Item(720,5).Value = "No"
$ws.Cells.Item(720,6).Value = 17832
$ws.Cells.Item(720,7).Value = "State Farm Arena"
$ws.Cells.Item(720,8).Value = "Atlanta Hawks"
$ws.Cells.Item(720,9).Value = "Phoenix Suns"

$ws.Cells.Item(721,1).Value = "Sacramento Kings"
$ws.Cells.Item(721,2).Value = 133
$ws.Cells.Item(721,2).NumberFormat = "#,##0"
$ws.Cells.Item(721,3).Value = "Indiana Pacers"
$ws.Cells.Item(721,4).Value = 122
$ws.Cells.Item(721,4).NumberFormat = "#,##0"
$ws.Cells.Item(721,5).Value = "No"
$ws.Cells.Item(721,6).Value = 17832
$ws.Cells.Item(721,7).Value = "Gainbridge Fieldhouse"
$ws.Cells.Item(721,8).Value = "Sacramento Kings"
$ws.Cells.Item(721,9).Value = "Indiana Pacers"

$ws.Cells.Item(722,1).Value = "Toronto Raptors"
$ws.Cells.Item(722,2).Value = 106
$ws.Cells.Item(722,2).NumberFormat = "#,##0"
$ws.Cells.Item(722,3).Value = "Houston Rockets"
$ws.Cells.Item(722,4).Value = 135
$ws.Cells.Item(722,4).NumberFormat = "#,##0"
$ws.Cells.Item(722,5).Value = "No"
$ws.Cells.Item(722,6).Value = 17832
$ws.Cells.Item(722,7).Value = "Toyota Center"
$ws.Cells.Item(722,8).Value = "Houston Rockets"
$ws.Cells.Item(722,9).Value = "Toronto Raptors"

$ws.Cells.Item(723,1).Value = "Golden State Warriors"
$ws.Cells.Item(723,2).Value = 121
$ws.Cells.Item(723,2).NumberFormat = "#,##0"
$ws.Cells.Item(723,3).Value = "Memphis Grizzlies"
$ws.Cells.Item(723,4).Value = 101
$ws.Cells.Item(723,4).NumberFormat = "#,##0"
$ws.Cells.Item(723,5).Value = "No"
$ws.Cells.Item(723,6).Value = 17832
$ws.Cells.Item(723,7).Value = "FedEx Forum"
$ws.Cells.Item(723,8).Value = "Golden State Warriors"
$ws.Cells.Item(723,9).Value = "Memphis Grizzlies"

$ws.Cells.Item(724,1).Value = "Orlando Magic"
$ws.Cells.Item(724,2).Value = 108
$ws.Cells.Item(724,2).NumberFormat = "#,##0"
$ws.Cells.Item(724,3).Value = "Minnesota Timberwolves"
$ws.Cells.Item(724,4).Value = 106
$ws.Cells.Item(724,4).NumberFormat = "#,##0"
$ws.Cells.Item(724,5).Value = "No"
$ws.Cells.Item(724,6).Value = 17832
$ws.Cells.Item(724,7).Value = "Target Center"
$ws.Cells.Item(724,8).Value = "Orlando Magic"
$ws.Cells.Item(724,9).Value = "Minnesota Timberwolves"

$ws.Cells.Item(725,1).Value = "Charlotte Hornets"
$ws.Cells.Item(725,2).Value = 106
$ws.Cells.Item(725,2).NumberFormat = "#,##0"
$ws.Cells.Item(725,3).Value = "Oklahoma City Thunder"
$ws.Cells.Item(725,4).Value = 126
$ws.Cells.Item(725,4).NumberFormat = "#,##0"
$ws.Cells.Item(725,5).Value = "No"
$ws.Cells.Item(725,6).Value = 17832
$ws.Cells.Item(725,7).Value = "Paycom Center"
$ws.Cells.Item(725,8).Value = "Oklahoma City Thunder"
$ws.Cells.Item(725,9).Value = "Charlotte Hornets"

$ws.Cells.Item(726,1).Value = "New Orleans Pelicans"
$ws.Cells.Item(726,2).Value = 114
$ws.Cells.Item(726,2).NumberFormat = "#,##0"
$ws.Cells.Item(726,3).Value = "San Antonio Spurs"
$ws.Cells.Item(726,4).Value = 113
$ws.Cells.Item(726,4).NumberFormat = "#,##0"
$ws.Cells.Item(726,5).Value = "No"
$ws.Cells.Item(726,6).Value = 17832
$ws.Cells.Item(726,7).Value = "Frost Bank Center"
$ws.Cells.Item(726,8).Value = "New Orleans Pelicans"
$ws.Cells.Item(726,9).Value = "San Antonio Spurs"

$ws.Cells.Item(727,1).Value = "Portland Trail Blazers"
$ws.Cells.Item(727,2).Value = 108
$ws.Cells.Item(727,2).NumberFormat = "#,##0"
$ws.Cells.Item(727,3).Value = "Denver Nuggets"
$ws.Cells.Item(727,4).Value = 120
$ws.Cells.Item(727,4).NumberFormat = "#,##0"
$ws.Cells.Item(727,5).Value = "No"
$ws.Cells.Item(727,6).Value = 17832
$ws.Cells.Item(727,7).Value = "Ball Arena"
$ws.Cells.Item(727,8).Value = "Denver Nuggets"
$ws.Cells.Item(727,9).Value = "Portland Trail Blazers"

# Leave the sheet scrolled/selected the same way the author left it.
$ws.Range("H719").Select() | Out-Null